# Working on salesperson reports code.
#
# The "Lookup" sheet holds, per canonical field (row 1 = canonical header),
# a list of known source-file column-name variants in the rows below it.
#
# Column N ("Split Percentage") previously also contained two commission-rate
# variants ("Comm %" and "Comm.%") plus a third ("Comm. %") further down -
# these belong with the Commission Rate synonyms in column O instead.
#
#   - "Comm %"   (N2) is a pure duplicate of other "Comm..." variants
#                 already listed under Commission Rate - just delete it.
#   - "Comm.%"   (N3) and "Comm. %" (N7) are moved down to the end of the
#                 Commission Rate list in column O (O12, O13).
#
# Deleting entries from column N shifts the remaining Split-Percentage
# synonyms up to close the gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookup")

# --- Column N: drop "Comm %", "Comm.%" and "Comm. %", shift the rest up ---
$ws.Range("N2").Value = "Split Percent"
$ws.Range("N3").Value = "Split%"
$ws.Range("N4").Value = "Split %"
$ws.Range("N5").Value = "Commission Percentage"
$ws.Range("N6").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("N8").ClearContents()

# --- Column O: append the two relocated commission-rate variants ---
$ws.Range("O12").Value = "Comm. %"
$ws.Range("O13").Value = "Comm.%"

# Reflect where the user was working (matches the scroll/selection state
# captured when the workbook was saved).
$ws.Range("N8").Select()
